$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions update).
# Several "Price" cells hold digit-grouped strings (e.g. "26.357.20") that
# are stored as text, not numbers. Where the refreshed price still parses as
# a plain decimal (e.g. "239.67"), force the cell to Text format first so
# Excel keeps storing our assignment as a literal string instead of silently
# converting it to a number.

$ws.Range("D2").Value = '26.356.77'
$ws.Range("E2").Value = '  +3.08%  '
$ws.Range("D3").Value = '1.718.08'
$ws.Range("E3").Value = '  +3.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.67'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4735'
$ws.Range("E7").Value = '  -1.20%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06212'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '1.715.76'
$ws.Range("E10").Value = '  +2.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07058'
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.34'
$ws.Range("E12").Value = '  +4.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5912'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.419'
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.04'
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '26.347.19'
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006821'
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").Value = '1.935.71'
$ws.Range("E21").Value = '  +2.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.543'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.761'
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.324'
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.46'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.404'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.03'
$ws.Range("E28").Value = '  +3.26%  '
$ws.Range("E29").Value = '  +4.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.997'
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.685'
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07750'
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04422'
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  +3.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6189'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9357'
$ws.Range("E37").Value = '  +9.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '112.89'
$ws.Range("E38").Value = '  +15.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.419'
$ws.Range("E39").Value = '  -7.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.918'
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.336'
$ws.Range("E43").Value = '  +13.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3810'
$ws.Range("E44").Value = '  +1.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1170'
$ws.Range("E45").Value = '  +4.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.300'
$ws.Range("E46").Value = '  +1.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05283'
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.29'
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.717'
$ws.Range("E49").Value = '  +5.02%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.214'
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3362'
$ws.Range("E51").Value = '  +1.03%  '
